# Weekly update for "Hortaliza, Vega Modelo de Temuco - Haba":
# a new price record (week) is inserted as row 14, pushing the existing
# rows 14-42 down to 15-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (shifts rows 14-42 -> 15-43,
# and extends the sheet dimension to A1:R43).
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44498
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 100112026
$ws.Range("G14").Value = "Haba"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 8000
$ws.Range("N14").Value = "$/saco 25 kilos"
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 320
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
